$wb = $excel.ActiveWorkbook

# "Generate Report for Handback": refresh the handoff/handback timestamps for the
# first (row 2) entry of each locale sheet.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 20:37:39"
$wsZhCn.Range("H2").Value = "2016-03-17 20:37:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 20:37:43"
$wsDeDe.Range("H2").Value = "2016-03-17 20:38:00"
